$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Address, $Text)
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$ws.Range("A2").Value = '[c]: ATP + Pantetheine4Phosphate ==> DPCOA + PPI'
Set-TextValue $ws "B2" '2884'
Set-TextValue $ws "C2" '0.00456666667'
$ws.Range("D2").Value = '{''numParticipants'': [2, 2], ''proximity'': 6.0, ''km'': '''', ''reactionID'': u''201'', ''entryID'': u''17911'', ''vmax'': u''0.00456666667'', ''reactionStoichiometry'': u''ATP + UMP = UDP + ADP'', ''species'': u''Streptococcus pneumoniae'', ''ECNumber'': u''2.7.4.22''}'
$ws.Range("E2").Value = '{''numParticipants'': [2, 2], ''proximity'': 6.0, ''km'': u''2.22E-5'', ''reactionID'': u''201'', ''entryID'': u''40434'', ''vmax'': u''5.33333333E-5'', ''reactionStoichiometry'': u''ATP + UMP = UDP + ADP'', ''species'': u''Streptococcus pneumoniae'', ''ECNumber'': u''2.7.4.22''}'
$ws.Range("F2").Value = '{''numParticipants'': [2, 2], ''proximity'': 6.0, ''km'': u''30.0'', ''reactionID'': u''201'', ''entryID'': u''51897'', ''vmax'': u''25.0'', ''reactionStoichiometry'': u''ATP + UMP = UDP + ADP'', ''species'': u''Bacillus subtilis'', ''ECNumber'': u''2.7.4.22''}'
Set-TextValue $ws "G2" '6.0'
$ws.Range("H2").Value = 'Lifted From 2.7.4'
$ws.Range("I2").Value = '11594, 11595, 11596, 17903, 17904, 17905, 17906, 17907, 17908, 17909, 17910, 17911, 17912, 17913, 17914, 17915, 17916, 17917, 17918, 17919, 17920, 17921, 17922, 17923, 17924, 17925, 17926, 17927, 17928, 17933, 17934, 17935, 40410, 40411, 40412, 40413, 40414, 40415, 40416, 40417, 40418, 40419, 40420, 40433, 40434, 40435, 40436, 49510, 51897, 51898, 51899'
$ws.Range("J2").Value = '0.001155, 0.00126666667, 0.001335, 0.00225, 0.00303333333, 0.00375, 0.00456666667, 0.0048, 0.00665, 0.00966666667, 0.0111, 0.6, 24.0, 25.0, 3.03333333E-4, 3.8E-4, 4.6E-4, 4.91666667E-4, 5.33333333E-5, 5.66666667E-4, 6.36666667E-4, 6.43333333E-4, 6.96666667E-4, 8.66666667E-4, 8.83333333E-5, 9.85E-4'
Set-TextValue $ws "K2" '1.0E-4'
$ws.Range("L2").Value = '{''numParticipants'': [2, 2], ''proximity'': 6.0, ''km'': u''1.0E-4'', ''reactionID'': u''201'', ''entryID'': u''17927'', ''vmax'': u''0.00665'', ''reactionStoichiometry'': u''ATP + UMP = UDP + ADP'', ''species'': u''Streptococcus pneumoniae'', ''ECNumber'': u''2.7.4.22''}'
$ws.Range("M2").Value = '{''numParticipants'': [2, 2], ''proximity'': 6.0, ''km'': u''1.0E-5'', ''reactionID'': u''201'', ''entryID'': u''40435'', ''vmax'': u''8.83333333E-5'', ''reactionStoichiometry'': u''ATP + UMP = UDP + ADP'', ''species'': u''Bacillus subtilis'', ''ECNumber'': u''2.7.4.22''}'
$ws.Range("N2").Value = '{''numParticipants'': [2, 2], ''proximity'': 6.0, ''km'': u''140.0'', ''reactionID'': u''13616'', ''entryID'': u''51899'', ''vmax'': u''0.6'', ''reactionStoichiometry'': u"ATP + 6-Azauridine 5''-phosphate = ADP + 6-Azauridine 5''-diphosphate", ''species'': u''Bacillus subtilis'', ''ECNumber'': u''2.7.4.22''}'
Set-TextValue $ws "O2" '6.0'
$ws.Range("P2").Value = 'Lifted From 2.7.4'
$ws.Range("Q2").Value = '11594, 11595, 11596, 17927, 17928, 40415, 40418, 40419, 40420, 40433, 40434, 40435, 40436, 49510, 51897, 51898, 51899, 51903, 51904'
$ws.Range("R2").Value = '0.9, 1.056E-4, 1.0E-4, 1.0E-5, 1.31E-4, 1.55E-4, 1.5E-4, 1.6E-5, 120.0, 140.0, 2.22E-5, 2.74E-5, 3.6E-5, 3.7E-5, 30.0, 7.6E-5, 8.0, 9.97E-5'
$ws.Range("A3").Value = '[c]: e3dCMP + H2O ==> e3dC + PI'
Set-TextValue $ws "C3" '4.33333333E-5'
$ws.Range("D3").Value = '{''numParticipants'': [2, 2], ''proximity'': 6.0, ''km'': u''9.0E-6'', ''reactionID'': u''581'', ''entryID'': u''22920'', ''vmax'': u''4.33333333E-5'', ''reactionStoichiometry'': u''D-Mannitol 1-phosphate + H2O = Phosphate + D-Mannitol'', ''species'': u''Streptococcus bovis'', ''ECNumber'': u''3.1.3.22''}'
$ws.Range("E3").Value = '{''numParticipants'': [2, 2], ''proximity'': 6.0, ''km'': u''0.07'', ''reactionID'': u''7713'', ''entryID'': u''22921'', ''vmax'': u''9.5E-6'', ''reactionStoichiometry'': u''H2O + Glucose 6-phosphate = Glucose + Phosphate'', ''species'': u''Streptococcus bovis'', ''ECNumber'': u''3.1.3.9''}'
$ws.Range("F3").Value = '{''numParticipants'': [2, 2], ''proximity'': 6.0, ''km'': u''0.011'', ''reactionID'': u''2493'', ''entryID'': u''16039'', ''vmax'': u''3.5E-4'', ''reactionStoichiometry'': u''2-Deoxy-D-glucose 6-phosphate + H2O = 2-Deoxy-D-glucose + Phosphate'', ''species'': u''Lactococcus lactis subsp. lactis'', ''ECNumber'': u''3.1.3.2''}'
Set-TextValue $ws "G3" '6.0'
$ws.Range("H3").Value = 'Lifted From 3.1.3'
$ws.Range("I3").Value = '16039, 16871, 16872, 22920, 22921, 22922, 42747'
$ws.Range("J3").Value = '2.91666667E-5, 3.33333333E-5, 3.5E-4, 4.33333333E-5, 6.33333333E-5, 9.0E-5, 9.5E-6'
Set-TextValue $ws "K3" '2.6E-6'
$ws.Range("L3").Value = '{''numParticipants'': [2, 2], ''proximity'': 1.0, ''km'': u''2.6E-6'', ''reactionID'': u''295'', ''entryID'': u''30219'', ''vmax'': u'''', ''reactionStoichiometry'': u''H2O + AMP = Phosphate + Adenosine'', ''species'': u''Mycoplasma fermentans'', ''ECNumber'': u''3.1.3.5''}'
$ws.Range("M3").Value = '{''numParticipants'': [2, 2], ''proximity'': 1.0, ''km'': u''6.3E-7'', ''reactionID'': u''295'', ''entryID'': u''30226'', ''vmax'': u'''', ''reactionStoichiometry'': u''H2O + AMP = Phosphate + Adenosine'', ''species'': u''Mycoplasma pulmonis'', ''ECNumber'': u''3.1.3.5''}'
$ws.Range("N3").Value = '{''numParticipants'': [2, 2], ''proximity'': 1.0, ''km'': u''1.04E-5'', ''reactionID'': u''295'', ''entryID'': u''30224'', ''vmax'': u'''', ''reactionStoichiometry'': u''H2O + AMP = Phosphate + Adenosine'', ''species'': u''Mycoplasma fermentans'', ''ECNumber'': u''3.1.3.5''}'
Set-TextValue $ws "O3" '1.0'
$ws.Range("P3").Value = 'Lifted From 3.1.3'
$ws.Range("Q3").Value = '30218, 30219, 30220, 30221, 30222, 30223, 30224, 30225, 30226, 30227'
$ws.Range("R3").Value = '1.04E-5, 1.0E-6, 1.8E-6, 2.2E-6, 2.6E-6, 4.1E-6, 4.2E-6, 6.3E-7'
